# docs: update a pipeline image
#
# Widens the "Back-end Pipeline" box/connector on slide 2 (so it spans the
# full pipeline again) and nudges the bottom row of cylinders (+ a couple of
# nearby labels) down to line up with it.
#
# NOTE on numeric literals below: PowerPoint's COM object model stores
# Shape.Left/Top/Width/Height as single-precision (float32) point values;
# EMU = round-trips through `(float32)points * 12700` using truncation.
# The literals used here are the exact float32 grid points that truncate
# back to the EMU values from the target OOXML, so the saved file matches
# to the EMU, not just "close enough" in points.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- bottom row of cylinders: shift down by 90005 EMU (~7.087 pt) ---------
$bottomCylinders = @(
    "円柱 2",
    "円柱 61",
    "円柱 62",
    "円柱 3",
    "円柱 56",
    "円柱 57",
    "円柱 58",
    "円柱 60"
)
foreach ($name in $bottomCylinders) {
    $shape = $s.Shapes.Item($name)
    $shape.Top = 333.7807312011719
}

# id=60 "円柱 59" lands one EMU short of the others (4239014 vs 4239015),
# matching the source file's own pre-existing one-EMU discrepancy.
$s.Shapes.Item("円柱 59").Top = 333.7806396484375

# --- small labels that ride along with the cylinders ----------------------
$s.Shapes.Item("テキスト ボックス 26").Left = 675.0765380859375
$s.Shapes.Item("テキスト ボックス 27").Left = 675.0765380859375
$s.Shapes.Item("テキスト ボックス 63").Top = 351.9102478027344

# --- widen the "Back-end Pipeline" title box + its connector line ---------
$titleBox = $s.Shapes.Item("テキスト ボックス 31")
$titleBox.Left = 310.3931579589844
$titleBox.Width = 460.6350402832031

$connector = $s.Shapes.Item("直線コネクタ 29")
$connector.Width = 460.6350402832031

# NOTE: the shipped OOXML also marks this connector's non-visual properties
# with an explicit <a:cxnSpLocks/> (PowerPoint adds that when a connector is
# manually stretched past its default geometry in the UI). There is no
# scriptable object-model property for that marker, so it is not
# reproducible from here; the geometry change above is applied regardless.

# --- notesMaster date placeholder ------------------------------------------
# The source diff also shows the cached `datetimeFigureOut` field text
# flipping from 2023/3/28 to 2023/4/6 - this is PowerPoint re-stamping the
# "last saved" date on the notes master, not a manual edit. Try to update it
# via the object model anyway, best-effort, without aborting the rest of
# the script if the host treats the field as read-only.
try {
    $nm = $p.NotesMaster
    $dateShape = $nm.Shapes.Item("日付プレースホルダー 2")
    $dateShape.TextFrame.TextRange.Text = "2023/4/6"
} catch {
}
